$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old standalone row that only held "8711623 - Denize Kalempa"
# in columns B/C (row 13). This shifts all following rows up by one,
# which realigns the label column (A) with its correct row per the
# target layout, while leaving the B/C "value" cells pointing at their
# old (now displaced) content - exactly matching the authored edit.
$ws.Range("A13").EntireRow.Delete()

# Now overwrite the specific value cells (col B and C) that differ from
# the simple "shift up by one" result, per the target workbook.

# Row 10: Objetivos: value becomes the docente code/name string.
$ws.Range("B10").Value = "8711623 - Denize Kalempa"
$ws.Range("C10").Value = "8711623 - Denize Kalempa"

# Row 13: Programa resumido: value becomes "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15: Programa: value becomes the activation date string.
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

# Row 18: Método: value becomes the docente code/name string.
$ws.Range("B18").Value = "8711623 - Denize Kalempa"
$ws.Range("C18").Value = "8711623 - Denize Kalempa"

# Row 19: Critério: value becomes the NF evaluation composition text.
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# Row 20: Norma de recuperação: value becomes "NF≥ 5,0."
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

# Row 21: Bibliografia: value becomes the recovery norm formula text.
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
